$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.297.55'
$ws.Range("E2").Value = '  -0.13%  '

# Row 3
$ws.Range("D3").Value = '1.692.40'
$ws.Range("E3").Value = '  +0.68%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.13%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5375'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.75%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.006'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2726'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.10%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06436'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.53%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.94%  '

# Row 11
$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '1.988.98'
$ws.Range("E11").Value = '  +18.28%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07694'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.44%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.529'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.24%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5804'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.20%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008379'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.61%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.27%  '

# Row 17
$ws.Range("D17").Value = '26.326.82'
$ws.Range("E17").Value = '  -0.09%  '

# Row 18
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.913'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.23%  '

# Row 19
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.03%  '

# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.03%  '

# Row 21
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.06%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.271'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.97%  '

# Row 23
$ws.Range("E23").Value = '  -0.10%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.71%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1291'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.67%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.868'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.09%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.50%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.389'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.87%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06149'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.66%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.326'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.01%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.602'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.30%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.588'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.13%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.690'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.74%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.033'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.41%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6202'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.30%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.425'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.69%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.761'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.84%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01646'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.45%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.174'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.90%  '

# Row 40
$ws.Range("D40").Value = '1.111.01'
$ws.Range("E40").Value = '  -0.45%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8781'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.42%  '

# Row 42
$ws.Range("E42").Value = '  -0.41%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.49%  '

# Row 44
$ws.Range("D44").Value = '1.842.29'
$ws.Range("E44").Value = '  +0.70%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.46%  '

# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000109'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.05%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.009'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.14%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.138'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.60%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05286'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.34%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4292'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.09%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.061'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.47%  '
